$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
  "D2" = "288.55"
  "E2" = "-0.59%"
  "E3" = "0.61%"
  "D4" = "4.920"
  "E4" = "-0.44%"
  "D5" = "0.07383"
  "E5" = "2.38%"
  "D6" = "2.215"
  "E6" = "24.88%"
  "D7" = "7.701"
  "E7" = "0.44%"
  "D8" = "3.731"
  "E8" = "0.21%"
  "D9" = "0.9080"
  "E9" = "1.42%"
  "D10" = "0.08712"
  "E10" = "12.70%"
  "D11" = "0.1691"
  "E11" = "1.91%"
  "D12" = "0.08191"
  "E12" = "1.69%"
  "D13" = "0.03109"
  "E13" = "2.16%"
  "E14" = "-0.72%"
  "D15" = "0.001499"
  "E15" = "-0.31%"
  "D16" = "0.005812"
  "E16" = "2.01%"
  "E17" = "0.44%"
  "D18" = "2.091"
  "E18" = "0.50%"
  "D19" = "0.3328"
  "E19" = "0.43%"
  "D20" = "0.1296"
  "E20" = "-0.71%"
  "D21" = "3.823"
  "E21" = "-5.45%"
  "E22" = "1.08%"
  "D23" = "0.04547"
  "E23" = "0.71%"
  "D24" = "0.001210"
  "E24" = "-0.29%"
  "D25" = "0.004144"
  "E25" = "3.35%"
  "E26" = "4.13%"
  "D39" = "0.01588"
  "E39" = "-0.82%"
  "D40" = "0.04469"
  "E40" = "1.64%"
  "D41" = "0.007321"
  "E41" = "-1.57%"
  "D42" = "0.009550"
  "E42" = "24.66%"
  "D43" = "0.1325"
  "E43" = "1.25%"
  "D44" = "0.002222"
  "E44" = "7.90%"
  "D45" = "0.008415"
  "E45" = "-8.55%"
  "D46" = "0.00006114"
  "E46" = "3.37%"
  "E47" = "0.08%"
  "D48" = "2.192"
  "E48" = "-2.38%"
  "E49" = "-33.27%"
  "E50" = "0.08%"
  "E51" = "0.08%"
}

foreach ($cell in $changes.Keys) {
  $rng = $ws.Range($cell)
  $rng.NumberFormat = "@"
  $rng.Value = $changes[$cell]
  $rng.Style = "Normal"
}

Write-Host "Done applying" $changes.Count "changes"
